$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mair")

# New working-hours entry: Programming / Create-Script + Inserts
# Copy the number formats from row 5 (date/time formatted cells) down into
# row 6 first so the new values inherit the correct [h]:mm:ss / date styles,
# matching the "Good"/"Bad" banding used throughout the table.
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)

$ws.Range("B6").Value = 41956
$ws.Range("C6").Value = "Programming"
$ws.Range("D6").Value = "Create-Script + Inserts"
$ws.Range("E6").Value = 0.041666666666666664
$ws.Range("F6").Value = 0.027777777777777776

# Match the author's final selection on the Mair sheet.
[void]$ws.Range("K12").Select()
